$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Ag. Director Abubakar Bello Afegbua"
$ws.Range("B3").Value = "19:26:53 GMT+0100 (West Africa Standard Time)"
$ws.Range("C3").Value = "Fri Jun 07 2024"
$ws.Range("D3").Value = "NPC Headquaters"

$ws.Range("A4").Value = "Ag. Director Abubakar Bello Afegbua"
$ws.Range("B4").Value = "19:26:54 GMT+0100 (West Africa Standard Time)"
$ws.Range("C4").Value = "Fri Jun 07 2024"
$ws.Range("D4").Value = "NPC Headquaters"
